$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 120 (pushing the existing rows 120-185 down
# to 121-186, one new weekly observation for "Provincia del Elquí").
$ws.Rows.Item(120).Insert()

$ws.Range("A120").Value = 8
$ws.Range("B120").Value = "Terminal La Palmera de La Serena"
$ws.Range("C120").Value = "Coquimbo"
$ws.Range("D120").Value = 44452
$ws.Range("E120").Value = 4
$ws.Range("F120").Value = 100114013
$ws.Range("G120").Value = "Zanahoria"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 720
$ws.Range("K120").Value = 5000
$ws.Range("L120").Value = 5500
$ws.Range("M120").Value = 5250
$ws.Range("N120").Value = "$/saco 20 kilos"
$ws.Range("O120").Value = "Provincia del Elquí"
$ws.Range("P120").Value = 262
$ws.Range("Q120").Value = 20
$ws.Range("R120").Value = "Hortaliza"
